# Update the "Metadata" worksheet to reflect the new IG publication metadata:
#  - Version bump 5.0.0 -> 6.0.0
#  - Date refresh
#  - Publisher now populated ("Alvearie Team")
#  - Duplicated "Contact" row replaced by a single "Jurisdiction" row
#
# (The "Elements" worksheet is untouched - its only differences in the
#  underlying OOXML come from shared-string re-indexing caused by the
#  sharedStrings table edits above, not from any real content change.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# The sheet previously had two identical "Contact" / "No display for
# ContactDetail" rows (rows 10 and 11). Remove the second one so the
# table goes back down to a single row for that slot; Excel will shift
# every following row up by one automatically.
$ws.Rows.Item(11).Delete()

# Version: 5.0.0 -> 6.0.0
$ws.Cells.Item(3, 2).Value = "6.0.0"

# Date: refreshed publish timestamp
$ws.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher: now populated
$ws.Cells.Item(9, 2).Value = "Alvearie Team"

# The remaining "Contact" row (now row 10) becomes "Jurisdiction"
$ws.Cells.Item(10, 1).Value = "Jurisdiction"
$ws.Cells.Item(10, 2).Value = "United States of America"
